$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

# Column A holds the date as plain text (matching the existing rows,
# which are stored as inline/shared strings rather than real dates), so
# force text formatting before assigning the value to stop Excel's
# automatic date-literal detection, then drop back to the sheet's normal
# (unformatted) style so no stray per-cell format sticks around.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/01/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1405966619741709
$ws.Cells.Item($row, 3).Value = 0.8594033380258291
